$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New handoff identifiers / timestamps that replace the previous handoff run.
# ---------------------------------------------------------------------------
$oldGuid = "2884f280-caf2-4e65-8a5f-da671b7c46a8"
$newGuid = "67f04f20-08ef-4063-b3fb-5928918bc9e8"

$newReadyDatetime   = "2016-09-05 11:28:15"
$zhHandoffDatetime  = "2016-09-05 11:28:09"
$zeroDatetime       = "0001-01-01 00:00:00"

$zhHandoffFile = "$newGuid.2c535e90786f0e209497d31499a56d7de120eacf.zh-cn.xlf"
$deHandoffFile = "$newGuid.2c535e90786f0e209497d31499a56d7de120eacf.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = "$newGuid.md"

$ovHl = $overview.Hyperlinks.Item(1)
$ovHl.TextToDisplay = "e2e\$newGuid.md"

$overview.Range("G2").Value = $newReadyDatetime

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zhHl = $zh.Hyperlinks.Item(1)
$zhHl.TextToDisplay = "$newGuid.md"

$zh.Range("G2").Value = $zhHandoffFile
$zh.Range("H2").Value = $zhHandoffDatetime

$zh.Range("I2").Value = ""
$zh.Range("I2").Style = "Normal"
$zh.Range("J2").Value = ""

$zh.Range("K2").Value = $zeroDatetime

$zh.Columns.Item(9).ColumnWidth = 18.6506053379604
$zh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$deHl = $de.Hyperlinks.Item(1)
$deHl.TextToDisplay = "$newGuid.md"

$de.Range("G2").Value = $deHandoffFile
$de.Range("H2").Value = $newReadyDatetime

$de.Range("I2").Value = ""
$de.Range("I2").Style = "Normal"
$de.Range("J2").Value = ""

$de.Range("K2").Value = $zeroDatetime

$de.Columns.Item(9).ColumnWidth = 18.6506053379604
$de.Columns.Item(10).ColumnWidth = 21.7054770333426

$wb.Save()
